$p = $ppt.ActivePresentation
try {
  $p.HasNotesMaster = -1
  Write-Host "Set HasNotesMaster ok, now:" $p.HasNotesMaster
} catch {
  Write-Host "ERROR: $_"
}
